$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("player_parameter")

# Swap the Name values of the two boss HP bar parameter blocks
$ws.Range("B25").Value = "BOSSHPBAR_BACK"
$ws.Range("B31").Value = "BOSSHPBAR_FRONT"

# Update the view state (scroll position + active selection)
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B32").Select()
